# Insert a new data row at row 316 (pushes existing rows 316-389 down to
# 317-390, extending the used range from A1:R389 to A1:R390) and populate
# the newly-inserted row with a new "Ajo" (garlic) price record for the
# Feria Lagunitas de Puerto Montt market.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 316..389 down one row, creating a blank row 316.
$ws.Rows(316).Insert()

# Fill the new row 316 with the record's data.
$ws.Cells.Item(316, 1).Value  = 4
$ws.Cells.Item(316, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(316, 3).Value  = "Los Lagos"
$ws.Cells.Item(316, 4).Value  = 44932
$ws.Cells.Item(316, 5).Value  = 10
$ws.Cells.Item(316, 6).Value  = 100112003
$ws.Cells.Item(316, 7).Value  = "Ajo"
$ws.Cells.Item(316, 8).Value  = "Chino"
$ws.Cells.Item(316, 9).Value  = "Primera"
$ws.Cells.Item(316, 10).Value = 240
$ws.Cells.Item(316, 11).Value = 18000
$ws.Cells.Item(316, 12).Value = 18000
$ws.Cells.Item(316, 13).Value = 18000
$ws.Cells.Item(316, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(316, 15).Value = "China"
$ws.Cells.Item(316, 16).Value = 1800
$ws.Cells.Item(316, 17).Value = 10
$ws.Cells.Item(316, 18).Value = "Hortaliza"
